$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 25/26: Toncoin and LidoDAOToken swap places
Set-TextValue $ws.Cells.Item(25, 2) "Toncoin"
Set-TextValue $ws.Cells.Item(25, 3) "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Cells.Item(26, 2) "LidoDAOToken"
Set-TextValue $ws.Cells.Item(26, 3) "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"

# Price (D) and Volume(1h) (E) updates
$priceVolume = @{
  2 = @("24.868.32", "  +0.05%  ")
  3 = @("1.706.71", "  -0.11%  ")
  4 = @("1.003", "  -0.23%  ")
  5 = @("317.85", "  -0.21%  ")
  6 = @($null, "  -0.40%  ")
  7 = @("0.4028", "  +1.41%  ")
  8 = @("0.4074", "  -1.08%  ")
  9 = @("1.482", "  -2.36%  ")
  10 = @("1.003", "  -0.21%  ")
  11 = @("53.68", "  +0.14%  ")
  12 = @("0.08830", "  -1.93%  ")
  13 = @("26.43", "  +8.19%  ")
  14 = @("7.491", "  -2.77%  ")
  15 = @("8.149", "  -1.05%  ")
  16 = @("0.00001360", "  -1.27%  ")
  17 = @("1.705.66", "  -0.97%  ")
  18 = @("97.09", "  -3.38%  ")
  19 = @("0.07174", "  +0.07%  ")
  20 = @("21.23", "  +5.06%  ")
  21 = @("7.288", "  -3.17%  ")
  22 = @("1.002", "  -0.69%  ")
  23 = @("14.38", "  -1.46%  ")
  24 = @("24.869.30", "  +0.03%  ")
  25 = @("2.327", "  -0.88%  ")
  26 = @("2.927", "  -5.38%  ")
  27 = @("23.24", "  +0.51%  ")
  28 = @("6.356", "  +21.49%  ")
  29 = @("166.98", "  -0.09%  ")
  30 = @("145.95", "  +4.54%  ")
  31 = @("8.424", "  -9.72%  ")
  32 = @("2.233", "  +13.50%  ")
  33 = @("1.894.27", "  -0.74%  ")
  34 = @("0.08842", "  -3.18%  ")
  35 = @("0.03199", "  +4.97%  ")
  36 = @("7.245", "  -8.79%  ")
  37 = @("1.034", "  -5.18%  ")
  38 = @("0.2879", "  +1.81%  ")
  39 = @("0.8453", "  +3.85%  ")
  40 = @("10.89", "  -2.38%  ")
  41 = @("0.09349", "  -0.16%  ")
  42 = @("14.18", "  -3.21%  ")
  43 = @("1.472", "  -0.98%  ")
  44 = @("17.57", "  +4.32%  ")
  45 = @("2.712", "  +1.81%  ")
  46 = @("0.7437", "  -0.05%  ")
  47 = @("4.243", "  -0.85%  ")
  48 = @("1.403", "  +3.28%  ")
  49 = @("0.9982", "  -0.40%  ")
  50 = @("141.80", "  +0.40%  ")
  51 = @("0.08358", "  +3.54%  ")
}

foreach ($row in $priceVolume.Keys) {
    $vals = $priceVolume[$row]
    if ($vals[0] -ne $null) { Set-TextValue $ws.Cells.Item($row, 4) $vals[0] }
    if ($vals[1] -ne $null) { Set-TextValue $ws.Cells.Item($row, 5) $vals[1] }
}
